$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to Text format so numeric-looking strings
# (e.g. "104.50", "0.000007530", "1.003") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated values exactly as they appear in the target workbook.
$ws.Range("D2").Value = "30.674.09"
$ws.Range("D3").Value = "1.918.97"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "239.43"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "0.4780"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("D8").Value = "0.2883"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "0.06709"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "18.82"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").Value = "104.50"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07726"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.915.50"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "5.259"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "0.6867"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "267.45"
$ws.Range("E16").Value = "  -5.85%  "
$ws.Range("D17").Value = "30.676.29"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "0.000007530"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "12.75"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").Value = "5.467"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "0.4553"
$ws.Range("E23").Value = "  -9.14%  "
$ws.Range("D24").Value = "6.351"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("D25").Value = "9.729"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "163.42"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").Value = "19.06"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("D28").Value = "2.113"
$ws.Range("E28").Value = "  -4.53%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "0.1021"
$ws.Range("E29").Value = "  -3.49%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.395"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "1.527"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "4.429"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("D33").Value = "4.245"
$ws.Range("E33").Value = "  -4.69%  "
$ws.Range("D34").Value = "0.04761"
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("D35").Value = "0.7361"
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("D36").Value = "1.126"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").Value = "2.718"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "0.01961"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").Value = "2.645"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "6.349"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").Value = "75.65"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "2.015"
$ws.Range("E43").Value = "  -4.70%  "
$ws.Range("D44").Value = "0.8681"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "106.54"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "0.4326"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "7.588"
$ws.Range("E48").Value = "  -7.00%  "
$ws.Range("D49").Value = "964.92"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "35.37"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1207"
$ws.Range("E51").Value = "  -4.16%  "
